$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values according to updated election results
$ws.Range("I2").Value = 386
$ws.Range("J2").Value = 1624
$ws.Range("K2").Value = 4
$ws.Range("L2").Value = 447
$ws.Range("N2").Value = 274
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 5
$ws.Range("R2").Value = 18
$ws.Range("S2").Value = 184
$ws.Range("T2").Value = 247
$ws.Range("U2").Value = 30
$ws.Range("V2").Value = 2486
$ws.Range("W2").Value = 0
$ws.Range("X2").Value = 2590
$ws.Range("Y2").Value = 4
$ws.Range("Z2").Value = 29
$ws.Range("AA2").Value = 14
